$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-looking numeric/percentage values need a leading apostrophe so Excel
# keeps them as literal text (matching the workbook's inlineStr cells) instead
# of auto-converting to a number; the style is then reset so no stray number
# format sticks to the cell.
function Set-TextValue($cell, $value) {
    $ws.Range($cell).Value = "'" + $value
    $ws.Range($cell).Style = "Normal"
}

Set-TextValue 'D2' '279.42'
Set-TextValue 'E2' '0.60%'
Set-TextValue 'D3' '27.47'
Set-TextValue 'E3' '0.91%'
Set-TextValue 'D4' '4.836'
Set-TextValue 'E4' '-0.85%'
Set-TextValue 'D5' '0.06388'
Set-TextValue 'E5' '-0.12%'
Set-TextValue 'D6' '7.052'
Set-TextValue 'E6' '1.03%'
Set-TextValue 'D7' '1.300'
Set-TextValue 'E7' '4.05%'
Set-TextValue 'D8' '0.8961'
Set-TextValue 'E8' '1.61%'
Set-TextValue 'D9' '0.1544'
Set-TextValue 'E9' '1.62%'
Set-TextValue 'D10' '0.06242'
Set-TextValue 'E10' '21.59%'
Set-TextValue 'D11' '0.07532'
Set-TextValue 'E11' '0.03%'
Set-TextValue 'D12' '0.02942'
Set-TextValue 'E12' '-0.11%'
Set-TextValue 'D13' '0.09003'
Set-TextValue 'E13' '-0.14%'
Set-TextValue 'D14' '0.001572'
Set-TextValue 'E14' '0.28%'
Set-TextValue 'D15' '0.0006415'
Set-TextValue 'E15' '0.17%'
Set-TextValue 'D16' '0.006023'
Set-TextValue 'E16' '1.37%'
Set-TextValue 'E17' '0.78%'
Set-TextValue 'D18' '3.324'
Set-TextValue 'E18' '0.14%'
Set-TextValue 'D19' '2.229'
Set-TextValue 'E19' '-1.90%'
Set-TextValue 'E21' '1.06%'
Set-TextValue 'D22' '3.907'
Set-TextValue 'E22' '0.03%'
Set-TextValue 'E23' '0.05%'
Set-TextValue 'D24' '0.1503'
Set-TextValue 'E24' '8.89%'
Set-TextValue 'E25' '0.19%'
Set-TextValue 'D26' '0.004279'
Set-TextValue 'E26' '10.35%'
Set-TextValue 'D28' '0.0001180'
Set-TextValue 'E28' '-1.74%'
Set-TextValue 'D29' '0.0001653'
Set-TextValue 'E29' '-14.62%'
Set-TextValue 'D40' '0.04066'
Set-TextValue 'E40' '-2.27%'
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue 'D41' '0.006684'
Set-TextValue 'E41' '-2.47%'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue 'D42' '0.1406'
Set-TextValue 'E42' '19.22%'
Set-TextValue 'D43' '0.002090'
Set-TextValue 'E43' '2.88%'
Set-TextValue 'D44' '0.01101'
Set-TextValue 'E44' '-2.12%'
Set-TextValue 'D45' '0.00005549'
Set-TextValue 'E45' '6.98%'
Set-TextValue 'D46' '1.628'
Set-TextValue 'E46' '9.85%'
Set-TextValue 'D47' '0.01849'
Set-TextValue 'E47' '-8.73%'
